# Add support for enclosure door open sensor
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New wiring row for the enclosure door sensor, inserted between the
# "Soil T/VWC probes" (row 14) and "Tmpr./RH probe" (row 17) groups,
# reusing the merged F-column layout pattern used elsewhere on the sheet.
#
# Values are written in an order that keeps the new shared-string indices
# aligned with the canonical workbook (Encl. door sensor, state detect,
# excitation, C8, 5V).
$ws.Range("F15").Value = "Encl. door sensor"
$ws.Range("G15").Value = "state detect"
$ws.Range("G16").Value = "excitation"
$ws.Range("I15").Value = "C8"
$ws.Range("I16").Value = "5V"
$ws.Range("H15").Value = "white"
$ws.Range("H16").Value = "white"

# Center the new F15:F16 description cell like its sibling groups (F4:F6,
# F7:F9, F10:F12, F13:F14) and merge it across the two rows.
$ws.Range("F15:F16").HorizontalAlignment = -4108
$ws.Range("F15:F16").VerticalAlignment = -4108
$ws.Range("F15:F16").WrapText = $false
$ws.Range("F15:F16").Merge()

# Match the saved selection recorded in the workbook.
$ws.Range("F17").Select()
